$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "neighbourhood"

$names = @(
    'Torre Baró',
    'Vallvidrera, el Tibidabo i les Planes',
    'la Vall d''Hebron',
    'Canyelles',
    'la Trinitat Vella',
    'la Trinitat Nova',
    'Pedralbes',
    'la Guineueta',
    'Ciutat Meridiana',
    'Sarrià',
    'les Roquetes',
    'el Poblenou',
    'Porta',
    'la Prosperitat',
    'Sant Genís dels Agudells',
    'la Sagrera',
    'Provençals del Poblenou',
    'Sant Gervasi - la Bonanova',
    'Horta',
    'la Font de la Guatlla',
    'Vilapicina i la Torre Llobeta',
    'la Marina de Port',
    'la Verneda i la Pau',
    'la Vila de Gràcia',
    'Baró de Viver',
    'la Dreta de l''Eixample',
    'Diagonal Mar i el Front Marítim del Poblenou',
    'Navas',
    'la Barceloneta',
    'la Font d''en Fargues',
    'el Camp d''en Grassot i Gràcia Nova',
    'Can Baró',
    'el Clot',
    'Sants - Badal',
    'la Sagrada Família',
    'l''Antiga Esquerra de l''Eixample',
    'la Marina del Prat Vermell',
    'la Vila Olímpica del Poblenou',
    'la Bordeta',
    'Vallcarca i els Penitents',
    'el Putxet i el Farró',
    'el Parc i la Llacuna del Poblenou',
    'el Baix Guinardó',
    'la Salut',
    'Sant Antoni',
    'Sants',
    'el Guinardó',
    'la Teixonera',
    'Sant Pere, Santa Caterina i la Ribera',
    'el Poble Sec',
    'Sant Martí de Provençals',
    'les Corts',
    'el Camp de l''Arpa del Clot',
    'el Fort Pienc',
    'la Nova Esquerra de l''Eixample',
    'el Bon Pastor',
    'Sant Gervasi - Galvany',
    'les Tres Torres',
    'Verdun',
    'el Congrés i els Indians',
    'el Besòs i el Maresme',
    'el Carmel',
    'el Barri Gòtic',
    'Hostafrancs',
    'el Coll',
    'el Raval',
    'Montbau',
    'Sant Andreu',
    'la Maternitat i Sant Ramon',
    'el Turó de la Peira',
    'Can Peguera'
)

$values = @(
    96.33333333333333,
    96.21052631578948,
    96,
    95.5,
    94.91666666666667,
    94.85714285714286,
    94.25,
    94,
    94,
    93.64,
    93.6086956521739,
    93.17487684729063,
    93,
    93,
    93,
    92.7910447761194,
    92.63366336633663,
    92.6,
    92.52941176470588,
    92.36521739130434,
    92.23809523809524,
    92.17142857142858,
    92.05454545454545,
    92.03729603729603,
    92,
    91.93487270574305,
    91.82993197278911,
    91.72463768115942,
    91.65116279069767,
    91.6,
    91.5036231884058,
    91.48648648648648,
    91.39516129032258,
    91.33557046979865,
    91.25802879291251,
    91.23094425483504,
    91.22222222222223,
    91.21935483870968,
    91.19469026548673,
    91.16666666666667,
    91.12903225806451,
    91.05504587155963,
    91.03508771929825,
    91.02439024390245,
    90.90725326991677,
    90.89057750759878,
    90.72619047619048,
    90.67857142857143,
    90.62121212121212,
    90.59371492704825,
    90.5,
    90.49197860962566,
    90.2720848056537,
    90.25389755011136,
    90.2315340909091,
    90.15384615384616,
    89.65313653136532,
    89.63636363636364,
    89.57894736842105,
    89.5,
    89.39772727272727,
    89.38947368421053,
    89.23138832997988,
    89.13168724279835,
    89.05,
    88.99848599545798,
    88.71428571428571,
    88.26666666666667,
    87.45161290322581,
    85.6470588235294,
    80
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}